$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '54.351.87'
$ws.Cells.Item(2, 5).Value = '  -2.69%  '
$ws.Cells.Item(3, 4).Value = '2.286.16'
$ws.Cells.Item(3, 5).Value = '  -2.60%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.19%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '495.10'
$ws.Cells.Item(5, 5).Value = '  -1.95%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '127.41'
$ws.Cells.Item(6, 5).Value = '  -1.73%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.527'
$ws.Cells.Item(8, 5).Value = '  -1.67%  '
$ws.Cells.Item(9, 4).Value = '2.284.19'
$ws.Cells.Item(9, 5).Value = '  -3.12%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0941'
$ws.Cells.Item(10, 5).Value = '  -3.12%  '
$ws.Cells.Item(11, 5).Value = '  +0.29%  '
$ws.Cells.Item(12, 5).Value = '  +0.33%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.63'
$ws.Cells.Item(13, 5).Value = '  -3.48%  '
$ws.Cells.Item(14, 4).Value = '2.688.87'
$ws.Cells.Item(14, 5).Value = '  -2.69%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '21.51'
$ws.Cells.Item(15, 5).Value = '  +0.23%  '
$ws.Cells.Item(16, 4).Value = '54.123.15'
$ws.Cells.Item(16, 5).Value = '  -3.02%  '
$ws.Cells.Item(17, 5).Value = '  -2.40%  '
$ws.Cells.Item(18, 4).Value = '2.325.41'
$ws.Cells.Item(18, 5).Value = '  -0.99%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '9.93'
$ws.Cells.Item(19, 5).Value = '  +0.25%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '4.05'
$ws.Cells.Item(20, 5).Value = '  +0.93%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '301.82'
$ws.Cells.Item(21, 5).Value = '  -2.77%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.40'
$ws.Cells.Item(22, 5).Value = '  +3.23%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  +0.24%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '5.36'
$ws.Cells.Item(24, 5).Value = '  -2.79%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '63.78'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.62%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.373'
$ws.Cells.Item(27, 5).Value = '  +0.83%  '
$ws.Cells.Item(28, 4).Value = '2.388.76'
$ws.Cells.Item(28, 5).Value = '  -2.74%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.149'
$ws.Cells.Item(29, 5).Value = '  +2.25%  '
$ws.Cells.Item(30, 5).Value = '  +0.12%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '166.18'
$ws.Cells.Item(31, 5).Value = '  -3.20%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.60'
$ws.Cells.Item(32, 5).Value = '  -2.72%  '
$ws.Cells.Item(33, 4).Value = '0.0₃0682'
$ws.Cells.Item(33, 5).Value = '  -2.90%  '
$ws.Cells.Item(34, 5).Value = '  +1.67%  '
$ws.Cells.Item(35, 5).Value = '  -0.02%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.998'
$ws.Cells.Item(36, 5).Value = '  +0.19%  '
$ws.Cells.Item(37, 5).Value = '  +0.67%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '17.56'
$ws.Cells.Item(38, 5).Value = '  -0.62%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.19'
$ws.Cells.Item(39, 5).Value = '  +1.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.871'
$ws.Cells.Item(40, 5).Value = '  +5.11%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.62'
$ws.Cells.Item(41, 5).Value = '  -0.54%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '35.46'
$ws.Cells.Item(42, 5).Value = '  -1.73%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.40'
$ws.Cells.Item(44, 5).Value = '  +0.63%  '
$ws.Cells.Item(45, 5).Value = '  +0.04%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '126.06'
$ws.Cells.Item(46, 5).Value = '  +0.15%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '4.84'
$ws.Cells.Item(47, 5).Value = '  -0.53%  '
$ws.Cells.Item(48, 5).Value = '  -0.62%  '
$ws.Cells.Item(49, 5).Value = '  -2.40%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '238.07'
$ws.Cells.Item(50, 5).Value = '  -0.44%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0479'
$ws.Cells.Item(51, 5).Value = '  +0.49%  '
